$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row to append (row 80), continuing the daily update series
$newRow = 80
$prevRow = 79

# Copy formatting from the previous row so the new row matches existing style
$ws.Range("A$prevRow:D$prevRow").Copy() | Out-Null
$ws.Range("A$newRow:D$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the values for the new row
$ws.Cells.Item($newRow, 1).Value = 46029
$ws.Cells.Item($newRow, 2).Value = 180
$ws.Cells.Item($newRow, 3).Value = 189
$ws.Cells.Item($newRow, 4).Value = 184

$wb.Save()
